# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns,
# and apply the re-ranking of BabyDogeCoin / XinFinNetwork / Mantle / Aave (rows 46-49).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.135.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.839.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.42%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.52%  "
$ws.Range("E6").Value = "  -2.63%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2994"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07458"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.23"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07644"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.837.42"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.030"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.99%  "
$ws.Range("E14").Value = "  -2.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "87.70"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.142"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -7.68%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.132.63"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008217"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.083.04"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "230.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.333"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9999"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1431"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.701"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.502"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.74%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.260"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.137"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.195"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05342"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7533"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.846"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.12%  "
$ws.Range("E36").Value = "  -2.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.688"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.315.58"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.71%  "
$ws.Range("E39").Value = "  -2.95%  "
$ws.Range("E40").Value = "  -1.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9438"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.028"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "104.63"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9992"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.982.99"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.61%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000123"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.40%  "
$ws.Range("B47").Value = "XinFinNetwork"
$ws.Range("C47").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.07786"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +19.84%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5180"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.39%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "64.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.773"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.417"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.88%  "
